# Append three new CB stat blocks (Charvarius Ward, J.T. Gray, Mike Ford)
# to the existing player-comparison table, continuing the established
# alternating row-fill pattern (yellow/green/yellow...).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the most recent 3-row block (Denzel Ward, style 3 /
# yellow fill) down onto the next block (rows 11-13) so the alternating
# shading continues correctly: 11-13 yellow, 14-16 green, 17-19 yellow.
$ws.Range("A5:F7").Copy()
$ws.Range("A11:F13").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A8:F10").Copy()
$ws.Range("A14:F16").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A11:F13").Copy()
$ws.Range("A17:F19").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = $false

# --- Row 11: Charvarius Ward / Group1 ---
$ws.Cells.Item(11,1).Value2 = "Charvarius Ward"
$ws.Cells.Item(11,2).Value2 = "Group1"
$ws.Cells.Item(11,3).Value2 = 8.666666666666666
$ws.Cells.Item(11,4).Value2 = 64
$ws.Cells.Item(11,5).Value2 = 47.33333333333334
$ws.Cells.Item(11,6).Value2 = 16.66666666666667

# --- Row 12: Charvarius Ward / Group2 ---
$ws.Cells.Item(12,1).Value2 = "Charvarius Ward"
$ws.Cells.Item(12,2).Value2 = "Group2"
$ws.Cells.Item(12,3).Value2 = 13.66666666666667
$ws.Cells.Item(12,4).Value2 = 71
$ws.Cells.Item(12,5).Value2 = 51
$ws.Cells.Item(12,6).Value2 = 20

# --- Row 13: Charvarius Ward / Difference ---
$ws.Cells.Item(13,1).Value2 = "Charvarius Ward"
$ws.Cells.Item(13,2).Value2 = "Difference"
$ws.Cells.Item(13,3).Value2 = 5
$ws.Cells.Item(13,4).Value2 = 7
$ws.Cells.Item(13,5).Value2 = 3.666666666666664
$ws.Cells.Item(13,6).Value2 = 3.333333333333332

# --- Row 14: J.T. Gray / Group1 ---
$ws.Cells.Item(14,1).Value2 = "J.T. Gray"
$ws.Cells.Item(14,2).Value2 = "Group1"
$ws.Cells.Item(14,3).Value2 = 0.3333333333333333
$ws.Cells.Item(14,4).Value2 = 15.66666666666667
$ws.Cells.Item(14,5).Value2 = 13
$ws.Cells.Item(14,6).Value2 = 2.666666666666667

# --- Row 15: J.T. Gray / Group2 ---
$ws.Cells.Item(15,1).Value2 = "J.T. Gray"
$ws.Cells.Item(15,2).Value2 = "Group2"
$ws.Cells.Item(15,3).Value2 = 0
$ws.Cells.Item(15,4).Value2 = 18
$ws.Cells.Item(15,5).Value2 = 9.333333333333334
$ws.Cells.Item(15,6).Value2 = 8.666666666666666

# --- Row 16: J.T. Gray / Difference ---
$ws.Cells.Item(16,1).Value2 = "J.T. Gray"
$ws.Cells.Item(16,2).Value2 = "Difference"
$ws.Cells.Item(16,3).Value2 = -0.3333333333333333
$ws.Cells.Item(16,4).Value2 = 2.333333333333334
$ws.Cells.Item(16,5).Value2 = -3.666666666666666
$ws.Cells.Item(16,6).Value2 = 6

# --- Row 17: Mike Ford / Group1 ---
$ws.Cells.Item(17,1).Value2 = "Mike Ford"
$ws.Cells.Item(17,2).Value2 = "Group1"
$ws.Cells.Item(17,3).Value2 = 1
$ws.Cells.Item(17,4).Value2 = 10
$ws.Cells.Item(17,5).Value2 = 8.666666666666666
$ws.Cells.Item(17,6).Value2 = 1.333333333333333

# --- Row 18: Mike Ford / Group2 ---
$ws.Cells.Item(18,1).Value2 = "Mike Ford"
$ws.Cells.Item(18,2).Value2 = "Group2"
$ws.Cells.Item(18,3).Value2 = 0.6666666666666666
$ws.Cells.Item(18,4).Value2 = 18.66666666666667
$ws.Cells.Item(18,5).Value2 = 13.33333333333333
$ws.Cells.Item(18,6).Value2 = 5.333333333333333

# --- Row 19: Mike Ford / Difference ---
$ws.Cells.Item(19,1).Value2 = "Mike Ford"
$ws.Cells.Item(19,2).Value2 = "Difference"
$ws.Cells.Item(19,3).Value2 = -0.3333333333333334
$ws.Cells.Item(19,4).Value2 = 8.666666666666668
$ws.Cells.Item(19,5).Value2 = 4.666666666666668
$ws.Cells.Item(19,6).Value2 = 4
